# Enhance SnapMap PAY-GRADE-ELEMENTS schema sheet: add new numeric / frequency /
# country fields and expand the description for numeric fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAY-GRADE-ELEMENTS")

$numericDesc = "Column name is just a sample and can be changed to match your intake form variable. For numeric variables, provide comma-separated values in format: min,max,target (e.g., '100,200,150' where min=100, max=200, target=150). Min and max values are required for numeric fields; target is optional."
$stringDesc = "Column name is just a sample and can be changed. Refer to documentation."

# --- Insert two new rows right after row 3 (for bonusPayFrequency / basePayFrequency). ---
# Inserting at row 4 twice pushes the old row 4 ("currency") down to row 6,
# and the new blank rows inherit formatting from the row above (row 3).
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# --- Insert six new rows after the (now shifted) "currency" row (row 6), for the
# re-added salaryAmount plus the additional new numeric/country fields (rows 7-12). ---
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# --- Row 3: salaryAmount -> joiningBonus, with the new longer description. ---
$ws.Range("A3").Value = "joiningBonus"
$ws.Range("B3").Value = $numericDesc

# --- Row 4 (new): bonusPayFrequency ---
$ws.Range("A4").Value = "bonusPayFrequency"
$ws.Range("B4").Value = $stringDesc
$ws.Range("C4").Value = "String"
$ws.Range("D4").Value = "Optional"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

# --- Row 5 (new): basePayFrequency ---
$ws.Range("A5").Value = "basePayFrequency"
$ws.Range("B5").Value = $stringDesc
$ws.Range("C5").Value = "String"
$ws.Range("D5").Value = "Optional"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""

# Row 6 ("currency") already carries its original values/format down from the insert.

# --- Row 7 (new): salaryAmount (re-added as a numeric field) ---
$ws.Range("A7").Value = "salaryAmount"
$ws.Range("B7").Value = $numericDesc
$ws.Range("C7").Value = "Numeric"
$ws.Range("D7").Value = "Optional"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

# --- Row 8 (new): variableBonus ---
$ws.Range("A8").Value = "variableBonus"
$ws.Range("B8").Value = $numericDesc
$ws.Range("C8").Value = "Numeric"
$ws.Range("D8").Value = "Optional"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

# --- Row 9 (new): stockOptionsAmount ---
$ws.Range("A9").Value = "stockOptionsAmount"
$ws.Range("B9").Value = $numericDesc
$ws.Range("C9").Value = "Numeric"
$ws.Range("D9").Value = "Optional"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

# --- Row 10 (new): relocationBonus ---
$ws.Range("A10").Value = "relocationBonus"
$ws.Range("B10").Value = $numericDesc
$ws.Range("C10").Value = "Numeric"
$ws.Range("D10").Value = "Optional"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# --- Row 11 (new): totalCompensation ---
$ws.Range("A11").Value = "totalCompensation"
$ws.Range("B11").Value = $numericDesc
$ws.Range("C11").Value = "Numeric"
$ws.Range("D11").Value = "Optional"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# --- Row 12 (new): countryCode ---
$ws.Range("A12").Value = "countryCode"
$ws.Range("B12").Value = $stringDesc
$ws.Range("C12").Value = "String"
$ws.Range("D12").Value = "Optional"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

Write-Host "PAY-GRADE-ELEMENTS updated: rows 1-12 (dimension A1:F12)."
